# Generate Report for Handback
#
# The localization-status report is regenerated: the entry for
# "23a2d031-0502-4253-8716-8488c85be147" moves from "Ready for handoff"
# to "Handed back: in sync with en-US" (it now has a real handback file +
# handback datetime instead of the 0001-01-01 00:00:00 placeholder), and
# the three report sheets are re-sorted alphabetically by source file name
# (05ba5c4b, 23a2d031, 44533c18, 50f396db) instead of creation order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ovRows = @(
    @("05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", "Handed back: in sync with en-US", "Handed back: in sync with en-US", "2016-03-24 02:37:32"),
    @("23a2d031-0502-4253-8716-8488c85be147.md", "Handed back: in sync with en-US", "Handed back: in sync with en-US", "2016-03-24 02:40:17"),
    @("44533c18-eeca-4e90-af39-4b8a6044c26d.md", "Handed back: in sync with en-US", "Handed back: in sync with en-US", "2016-03-24 02:35:52"),
    @("50f396db-a00a-4d9f-bc95-ae77014c1455.md", "Handed back: in sync with en-US", "Handed back: in sync with en-US", "2016-03-24 02:37:32")
)

for ($i = 0; $i -lt $ovRows.Length; $i++) {
    $r = $i + 2
    $row = $ovRows[$i]
    $ov.Cells.Item($r, 1).Value = $row[0]
    $ov.Cells.Item($r, 2).Value = $row[1]
    $ov.Cells.Item($r, 3).Value = $row[2]
    $ov.Cells.Item($r, 4).Value = $row[3]
}

foreach ($hl in $ov.Hyperlinks) {
    $row = $hl.Range.Row
    $hl.TextToDisplay = $ovRows[$row - 2][0]
}

# ---------------------------------------------------------------------
# Sheets 2 & 3: zh-cn / de-de detail reports
# ---------------------------------------------------------------------
# Columns: A Source File Name, B File Extension, C Status,
#          D Latest Handoff File, E Latest Handoff Datetime,
#          F Latest Target File, G Latest Handback File,
#          H Latest Handback DateTime, I Reference Tokens,
#          J Handoff Reason, K Dependency From, L Error Detail

$zhRows = @(
    @("05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", ".md", "Handed back: in sync with en-US", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.zh-cn.xlf", "2016-03-24 02:37:23", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.zh-cn.xlf", "2016-03-24 02:38:20"),
    @("23a2d031-0502-4253-8716-8488c85be147.md", ".md", "Handed back: in sync with en-US", "23a2d031-0502-4253-8716-8488c85be147.afef86aa67eb67354f4348aa078ecf89540ad640.zh-cn.xlf", "2016-03-24 02:40:08", "23a2d031-0502-4253-8716-8488c85be147.md", "23a2d031-0502-4253-8716-8488c85be147.afef86aa67eb67354f4348aa078ecf89540ad640.zh-cn.xlf", "2016-03-24 02:42:14"),
    @("44533c18-eeca-4e90-af39-4b8a6044c26d.md", ".md", "Handed back: in sync with en-US", "44533c18-eeca-4e90-af39-4b8a6044c26d.f92687e118b1b8bcedf9bf55f125c343d031c11a.zh-cn.xlf", "2016-03-24 02:35:43", "44533c18-eeca-4e90-af39-4b8a6044c26d.md", "44533c18-eeca-4e90-af39-4b8a6044c26d.f92687e118b1b8bcedf9bf55f125c343d031c11a.zh-cn.xlf", "2016-03-24 02:36:23"),
    @("50f396db-a00a-4d9f-bc95-ae77014c1455.md", ".md", "Handed back: in sync with en-US", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.zh-cn.xlf", "2016-03-24 02:37:23", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.zh-cn.xlf", "2016-03-24 02:38:20")
)

$deRows = @(
    @("05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", ".md", "Handed back: in sync with en-US", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.de-de.xlf", "2016-03-24 02:37:32", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.de-de.xlf", "2016-03-24 02:38:36"),
    @("23a2d031-0502-4253-8716-8488c85be147.md", ".md", "Handed back: in sync with en-US", "23a2d031-0502-4253-8716-8488c85be147.afef86aa67eb67354f4348aa078ecf89540ad640.de-de.xlf", "2016-03-24 02:40:17", "23a2d031-0502-4253-8716-8488c85be147.md", "23a2d031-0502-4253-8716-8488c85be147.afef86aa67eb67354f4348aa078ecf89540ad640.de-de.xlf", "2016-03-24 02:42:28"),
    @("44533c18-eeca-4e90-af39-4b8a6044c26d.md", ".md", "Handed back: in sync with en-US", "44533c18-eeca-4e90-af39-4b8a6044c26d.f92687e118b1b8bcedf9bf55f125c343d031c11a.de-de.xlf", "2016-03-24 02:35:52", "44533c18-eeca-4e90-af39-4b8a6044c26d.md", "44533c18-eeca-4e90-af39-4b8a6044c26d.f92687e118b1b8bcedf9bf55f125c343d031c11a.de-de.xlf", "2016-03-24 02:36:37"),
    @("50f396db-a00a-4d9f-bc95-ae77014c1455.md", ".md", "Handed back: in sync with en-US", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.de-de.xlf", "2016-03-24 02:37:32", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", "05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.de-de.xlf", "2016-03-24 02:38:36")
)

function Update-DetailSheet($ws, $rows) {
    for ($i = 0; $i -lt $rows.Length; $i++) {
        $r = $i + 2
        $row = $rows[$i]
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $ws.Cells.Item($r, 7).Value = $row[6]
        $ws.Cells.Item($r, 8).Value = $row[7]
    }

    # Existing hyperlinks (columns A, D, F, G on rows 2-4 before this edit;
    # row 5 only had A and D) follow their row to the re-sorted data.
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        $row = $hl.Range.Row
        $col = $hl.Range.Column
        $data = $rows[$row - 2]
        if ($col -eq 1) {
            $hl.TextToDisplay = $data[0]
        } elseif ($col -eq 4) {
            $hl.TextToDisplay = $data[3]
        } elseif ($col -eq 6) {
            $hl.TextToDisplay = $data[5]
        } elseif ($col -eq 7) {
            $hl.TextToDisplay = $data[6]
        }
    }

    # Row 5 (now "50f396db...") gains Latest Target File / Latest Handback
    # File hyperlinks it didn't have before (it duplicates row 2's links,
    # matching the source report generator's behaviour).
    $row5 = $rows[3]
    $fCell = $ws.Cells.Item(5, 6)
    $gCell = $ws.Cells.Item(5, 7)
    $ws.Hyperlinks.Add($fCell, "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/5e9bd7161628046f9b644ce5cc1d294018d9a377/e2e/05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.md", "", "", $row5[5]) | Out-Null
    $ws.Hyperlinks.Add($gCell, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/5c06aa61af1152b408b66cf31c6f85229a3a2050/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/05ba5c4b-a2dc-4c05-8a2b-0e901d070f3d.aa3e0c61f5dfd93a4ec30865223d026ce20b8e62.zh-cn.xlf", "", "", $row5[6]) | Out-Null
}

$zhWs = $wb.Worksheets.Item("zh-cn")
Update-DetailSheet $zhWs $zhRows

$deWs = $wb.Worksheets.Item("de-de")
Update-DetailSheet $deWs $deRows
